$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr8 = New-Object 'object[,]' 1,10
$arr8[0,0] = "فصل سوم منتهی به 1399/09"
$arr8[0,1] = "فصل چهارم منتهی به 1399/12"
$arr8[0,2] = "فصل اول منتهی به 1400/03"
$arr8[0,3] = "فصل دوم منتهی به 1400/06"
$arr8[0,4] = "فصل سوم منتهی به 1400/09"
$arr8[0,5] = "فصل چهارم منتهی به 1400/12"
$arr8[0,6] = "فصل اول منتهی به 1401/03"
$arr8[0,7] = "فصل دوم منتهی به 1401/06"
$arr8[0,8] = "فصل سوم منتهی به 1401/09"
$arr8[0,9] = "فصل چهارم منتهی به 1401/12"
$ws.Range("E8:N8").Value = $arr8

$arr10 = New-Object 'object[,]' 1,10
$arr10[0,0] = "-"
$arr10[0,1] = "-"
$arr10[0,2] = "-"
$arr10[0,3] = "-"
$arr10[0,4] = "-"
$arr10[0,5] = 38501
$arr10[0,6] = 38559
$arr10[0,7] = 36501
$arr10[0,8] = 29536
$arr10[0,9] = 17600
$ws.Range("E10:N10").Value = $arr10

$arr11 = New-Object 'object[,]' 1,10
$arr11[0,0] = 25310
$arr11[0,1] = 104077
$arr11[0,2] = "-"
$arr11[0,3] = "-"
$arr11[0,4] = 39785
$arr11[0,5] = "-"
$arr11[0,6] = "-"
$arr11[0,7] = "-"
$arr11[0,8] = "-"
$arr11[0,9] = "-"
$ws.Range("E11:N11").Value = $arr11

$arr12 = New-Object 'object[,]' 1,10
$arr12[0,0] = "-"
$arr12[0,1] = "-"
$arr12[0,2] = "-"
$arr12[0,3] = 61684
$arr12[0,4] = "-"
$arr12[0,5] = "-"
$arr12[0,6] = "-"
$arr12[0,7] = "-"
$arr12[0,8] = "-"
$arr12[0,9] = "-"
$ws.Range("E12:N12").Value = $arr12

$arr13 = New-Object 'object[,]' 1,10
$arr13[0,0] = "-"
$arr13[0,1] = "-"
$arr13[0,2] = "-"
$arr13[0,3] = 0
$arr13[0,4] = 0
$arr13[0,5] = "-"
$arr13[0,6] = 0
$arr13[0,7] = 0
$arr13[0,8] = 0
$arr13[0,9] = 0
$ws.Range("E13:N13").Value = $arr13

$arr14 = New-Object 'object[,]' 1,10
$arr14[0,0] = "-"
$arr14[0,1] = "-"
$arr14[0,2] = "-"
$arr14[0,3] = "-"
$arr14[0,4] = "-"
$arr14[0,5] = 0
$arr14[0,6] = "-"
$arr14[0,7] = "-"
$arr14[0,8] = "-"
$arr14[0,9] = "-"
$ws.Range("E14:N14").Value = $arr14

$arr15 = New-Object 'object[,]' 1,10
$arr15[0,0] = 879367
$arr15[0,1] = "-"
$arr15[0,2] = "-"
$arr15[0,3] = 1656850
$arr15[0,4] = "-"
$arr15[0,5] = "-"
$arr15[0,6] = "-"
$arr15[0,7] = "-"
$arr15[0,8] = "-"
$arr15[0,9] = "-"
$ws.Range("E15:N15").Value = $arr15

$arr16 = New-Object 'object[,]' 1,10
$arr16[0,0] = "-"
$arr16[0,1] = "-"
$arr16[0,2] = "-"
$arr16[0,3] = "-"
$arr16[0,4] = 920959
$arr16[0,5] = 889585
$arr16[0,6] = "-"
$arr16[0,7] = "-"
$arr16[0,8] = "-"
$arr16[0,9] = "-"
$ws.Range("E16:N16").Value = $arr16

$arr17 = New-Object 'object[,]' 1,10
$arr17[0,0] = "-"
$arr17[0,1] = "-"
$arr17[0,2] = "-"
$arr17[0,3] = "-"
$arr17[0,4] = "-"
$arr17[0,5] = 245
$arr17[0,6] = 290
$arr17[0,7] = 59
$arr17[0,8] = 58
$arr17[0,9] = 193
$ws.Range("E17:N17").Value = $arr17

$arr18 = New-Object 'object[,]' 1,10
$arr18[0,0] = 173
$arr18[0,1] = 876
$arr18[0,2] = 144
$arr18[0,3] = "-"
$arr18[0,4] = 392
$arr18[0,5] = "-"
$arr18[0,6] = "-"
$arr18[0,7] = "-"
$arr18[0,8] = "-"
$arr18[0,9] = "-"
$ws.Range("E18:N18").Value = $arr18

$arr19 = New-Object 'object[,]' 1,10
$arr19[0,0] = "-"
$arr19[0,1] = "-"
$arr19[0,2] = "-"
$arr19[0,3] = 327
$arr19[0,4] = "-"
$arr19[0,5] = "-"
$arr19[0,6] = "-"
$arr19[0,7] = "-"
$arr19[0,8] = "-"
$arr19[0,9] = "-"
$ws.Range("E19:N19").Value = $arr19

$arr20 = New-Object 'object[,]' 1,10
$arr20[0,0] = "-"
$arr20[0,1] = 3246123
$arr20[0,2] = 862790
$arr20[0,3] = "-"
$arr20[0,4] = "-"
$arr20[0,5] = "-"
$arr20[0,6] = "-"
$arr20[0,7] = "-"
$arr20[0,8] = "-"
$arr20[0,9] = "-"
$ws.Range("E20:N20").Value = $arr20

$arr21 = New-Object 'object[,]' 1,10
$arr21[0,0] = "-"
$arr21[0,1] = "-"
$arr21[0,2] = "-"
$arr21[0,3] = 0
$arr21[0,4] = "-"
$arr21[0,5] = 0
$arr21[0,6] = 0
$arr21[0,7] = 0
$arr21[0,8] = 0
$arr21[0,9] = 0
$ws.Range("E21:N21").Value = $arr21

$arr22 = New-Object 'object[,]' 1,10
$arr22[0,0] = "-"
$arr22[0,1] = 9752
$arr22[0,2] = "-"
$arr22[0,3] = "-"
$arr22[0,4] = "-"
$arr22[0,5] = "-"
$arr22[0,6] = "-"
$arr22[0,7] = "-"
$arr22[0,8] = "-"
$arr22[0,9] = "-"
$ws.Range("E22:N22").Value = $arr22

$arr23 = New-Object 'object[,]' 1,10
$arr23[0,0] = "-"
$arr23[0,1] = "-"
$arr23[0,2] = "-"
$arr23[0,3] = "-"
$arr23[0,4] = "-"
$arr23[0,5] = 2498
$arr23[0,6] = 2576
$arr23[0,7] = 2132
$arr23[0,8] = 3379
$arr23[0,9] = 2470
$ws.Range("E23:N23").Value = $arr23

$arr24 = New-Object 'object[,]' 1,10
$arr24[0,0] = 2318
$arr24[0,1] = "-"
$arr24[0,2] = "-"
$arr24[0,3] = 4713
$arr24[0,4] = "-"
$arr24[0,5] = "-"
$arr24[0,6] = "-"
$arr24[0,7] = "-"
$arr24[0,8] = "-"
$arr24[0,9] = "-"
$ws.Range("E24:N24").Value = $arr24

$arr25 = New-Object 'object[,]' 1,10
$arr25[0,0] = "-"
$arr25[0,1] = "-"
$arr25[0,2] = "-"
$arr25[0,3] = "-"
$arr25[0,4] = 2591
$arr25[0,5] = "-"
$arr25[0,6] = "-"
$arr25[0,7] = "-"
$arr25[0,8] = "-"
$arr25[0,9] = "-"
$ws.Range("E25:N25").Value = $arr25

$arr26 = New-Object 'object[,]' 1,10
$arr26[0,0] = "-"
$arr26[0,1] = "-"
$arr26[0,2] = "-"
$arr26[0,3] = "-"
$arr26[0,4] = "-"
$arr26[0,5] = "-"
$arr26[0,6] = 879094
$arr26[0,7] = 842781
$arr26[0,8] = 1173827
$arr26[0,9] = 1036086
$ws.Range("E26:N26").Value = $arr26

$arr27 = New-Object 'object[,]' 1,10
$arr27[0,0] = 907168
$arr27[0,1] = 3360828
$arr27[0,2] = 862934
$arr27[0,3] = 1723574
$arr27[0,4] = 963727
$arr27[0,5] = 930829
$arr27[0,6] = 920519
$arr27[0,7] = 881473
$arr27[0,8] = 1206800
$arr27[0,9] = 1056349
$ws.Range("E27:N27").Value = $arr27

$arr31 = New-Object 'object[,]' 1,10
$arr31[0,0] = "فصل سوم منتهی به 1399/09"
$arr31[0,1] = "فصل چهارم منتهی به 1399/12"
$arr31[0,2] = "فصل اول منتهی به 1400/03"
$arr31[0,3] = "فصل دوم منتهی به 1400/06"
$arr31[0,4] = "فصل سوم منتهی به 1400/09"
$arr31[0,5] = "فصل چهارم منتهی به 1400/12"
$arr31[0,6] = "فصل اول منتهی به 1401/03"
$arr31[0,7] = "فصل دوم منتهی به 1401/06"
$arr31[0,8] = "فصل سوم منتهی به 1401/09"
$arr31[0,9] = "فصل چهارم منتهی به 1401/12"
$ws.Range("E31:N31").Value = $arr31

$arr33 = New-Object 'object[,]' 1,10
$arr33[0,0] = 24046
$arr33[0,1] = 25132
$arr33[0,2] = 23866
$arr33[0,3] = "-"
$arr33[0,4] = "-"
$arr33[0,5] = "-"
$arr33[0,6] = 34913
$arr33[0,7] = 37300
$arr33[0,8] = 28762
$arr33[0,9] = 14766
$ws.Range("E33:N33").Value = $arr33

$arr34 = New-Object 'object[,]' 1,10
$arr34[0,0] = "-"
$arr34[0,1] = "-"
$arr34[0,2] = "-"
$arr34[0,3] = 36995
$arr34[0,4] = "-"
$arr34[0,5] = "-"
$arr34[0,6] = "-"
$arr34[0,7] = "-"
$arr34[0,8] = "-"
$arr34[0,9] = "-"
$ws.Range("E34:N34").Value = $arr34

$arr35 = New-Object 'object[,]' 1,10
$arr35[0,0] = "-"
$arr35[0,1] = "-"
$arr35[0,2] = "-"
$arr35[0,3] = "-"
$arr35[0,4] = 38340
$arr35[0,5] = "-"
$arr35[0,6] = "-"
$arr35[0,7] = "-"
$arr35[0,8] = "-"
$arr35[0,9] = "-"
$ws.Range("E35:N35").Value = $arr35

$arr36 = New-Object 'object[,]' 1,10
$arr36[0,0] = "-"
$arr36[0,1] = "-"
$arr36[0,2] = "-"
$arr36[0,3] = 0
$arr36[0,4] = 0
$arr36[0,5] = "-"
$arr36[0,6] = 0
$arr36[0,7] = 0
$arr36[0,8] = 0
$arr36[0,9] = 0
$ws.Range("E36:N36").Value = $arr36

$arr37 = New-Object 'object[,]' 1,10
$arr37[0,0] = 0
$arr37[0,1] = 0
$arr37[0,2] = "-"
$arr37[0,3] = "-"
$arr37[0,4] = "-"
$arr37[0,5] = "-"
$arr37[0,6] = "-"
$arr37[0,7] = "-"
$arr37[0,8] = "-"
$arr37[0,9] = "-"
$ws.Range("E37:N37").Value = $arr37

$arr38 = New-Object 'object[,]' 1,10
$arr38[0,0] = 876301
$arr38[0,1] = 928853
$arr38[0,2] = 854844
$arr38[0,3] = "-"
$arr38[0,4] = "-"
$arr38[0,5] = 893410
$arr38[0,6] = "-"
$arr38[0,7] = "-"
$arr38[0,8] = "-"
$arr38[0,9] = "-"
$ws.Range("E38:N38").Value = $arr38

$arr39 = New-Object 'object[,]' 1,10
$arr39[0,0] = "-"
$arr39[0,1] = "-"
$arr39[0,2] = "-"
$arr39[0,3] = 838244
$arr39[0,4] = "-"
$arr39[0,5] = "-"
$arr39[0,6] = "-"
$arr39[0,7] = "-"
$arr39[0,8] = "-"
$arr39[0,9] = "-"
$ws.Range("E39:N39").Value = $arr39

$arr40 = New-Object 'object[,]' 1,10
$arr40[0,0] = "-"
$arr40[0,1] = "-"
$arr40[0,2] = "-"
$arr40[0,3] = "-"
$arr40[0,4] = 882233
$arr40[0,5] = "-"
$arr40[0,6] = "-"
$arr40[0,7] = "-"
$arr40[0,8] = "-"
$arr40[0,9] = "-"
$ws.Range("E40:N40").Value = $arr40

$arr41 = New-Object 'object[,]' 1,10
$arr41[0,0] = 200
$arr41[0,1] = 179
$arr41[0,2] = 144
$arr41[0,3] = "-"
$arr41[0,4] = "-"
$arr41[0,5] = "-"
$arr41[0,6] = 355
$arr41[0,7] = 59
$arr41[0,8] = 58
$arr41[0,9] = 146
$ws.Range("E41:N41").Value = $arr41

$arr42 = New-Object 'object[,]' 1,10
$arr42[0,0] = "-"
$arr42[0,1] = "-"
$arr42[0,2] = "-"
$arr42[0,3] = 362
$arr42[0,4] = "-"
$arr42[0,5] = "-"
$arr42[0,6] = "-"
$arr42[0,7] = "-"
$arr42[0,8] = "-"
$arr42[0,9] = "-"
$ws.Range("E42:N42").Value = $arr42

$arr43 = New-Object 'object[,]' 1,10
$arr43[0,0] = "-"
$arr43[0,1] = "-"
$arr43[0,2] = "-"
$arr43[0,3] = "-"
$arr43[0,4] = 386
$arr43[0,5] = "-"
$arr43[0,6] = "-"
$arr43[0,7] = "-"
$arr43[0,8] = "-"
$arr43[0,9] = "-"
$ws.Range("E43:N43").Value = $arr43

$arr44 = New-Object 'object[,]' 1,10
$arr44[0,0] = 0
$arr44[0,1] = 0
$arr44[0,2] = "-"
$arr44[0,3] = "-"
$arr44[0,4] = "-"
$arr44[0,5] = "-"
$arr44[0,6] = "-"
$arr44[0,7] = "-"
$arr44[0,8] = "-"
$arr44[0,9] = "-"
$ws.Range("E44:N44").Value = $arr44

$arr46 = New-Object 'object[,]' 1,10
$arr46[0,0] = 0
$arr46[0,1] = 0
$arr46[0,2] = "-"
$arr46[0,3] = 0
$arr46[0,4] = "-"
$arr46[0,5] = "-"
$arr46[0,6] = 0
$arr46[0,7] = 0
$arr46[0,8] = 0
$arr46[0,9] = 0
$ws.Range("E46:N46").Value = $arr46

$arr49 = New-Object 'object[,]' 1,10
$arr49[0,0] = "-"
$arr49[0,1] = "-"
$arr49[0,2] = "-"
$arr49[0,3] = "-"
$arr49[0,4] = "-"
$arr49[0,5] = "-"
$arr49[0,6] = 2136
$arr49[0,7] = 1568
$arr49[0,8] = 3451
$arr49[0,9] = 2175
$ws.Range("E49:N49").Value = $arr49

$arr50 = New-Object 'object[,]' 1,10
$arr50[0,0] = 2497
$arr50[0,1] = 2112
$arr50[0,2] = 2614
$arr50[0,3] = "-"
$arr50[0,4] = "-"
$arr50[0,5] = "-"
$arr50[0,6] = "-"
$arr50[0,7] = "-"
$arr50[0,8] = "-"
$arr50[0,9] = "-"
$ws.Range("E50:N50").Value = $arr50

$arr51 = New-Object 'object[,]' 1,10
$arr51[0,0] = "-"
$arr51[0,1] = "-"
$arr51[0,2] = "-"
$arr51[0,3] = 2161
$arr51[0,4] = "-"
$arr51[0,5] = "-"
$arr51[0,6] = "-"
$arr51[0,7] = "-"
$arr51[0,8] = "-"
$arr51[0,9] = "-"
$ws.Range("E51:N51").Value = $arr51

$arr52 = New-Object 'object[,]' 1,10
$arr52[0,0] = "-"
$arr52[0,1] = "-"
$arr52[0,2] = "-"
$arr52[0,3] = "-"
$arr52[0,4] = 2401
$arr52[0,5] = "-"
$arr52[0,6] = "-"
$arr52[0,7] = "-"
$arr52[0,8] = "-"
$arr52[0,9] = "-"
$ws.Range("E52:N52").Value = $arr52

$arr53 = New-Object 'object[,]' 1,10
$arr53[0,0] = "-"
$arr53[0,1] = "-"
$arr53[0,2] = "-"
$arr53[0,3] = "-"
$arr53[0,4] = "-"
$arr53[0,5] = "-"
$arr53[0,6] = 863861
$arr53[0,7] = 758681
$arr53[0,8] = 1177714
$arr53[0,9] = 819505
$ws.Range("E53:N53").Value = $arr53

$arr54 = New-Object 'object[,]' 1,10
$arr54[0,0] = 903044
$arr54[0,1] = 956276
$arr54[0,2] = 881468
$arr54[0,3] = 877762
$arr54[0,4] = 923360
$arr54[0,5] = 893410
$arr54[0,6] = 901265
$arr54[0,7] = 797608
$arr54[0,8] = 1209985
$arr54[0,9] = 836592
$ws.Range("E54:N54").Value = $arr54

$arr58 = New-Object 'object[,]' 1,10
$arr58[0,0] = "فصل سوم منتهی به 1399/09"
$arr58[0,1] = "فصل چهارم منتهی به 1399/12"
$arr58[0,2] = "فصل اول منتهی به 1400/03"
$arr58[0,3] = "فصل دوم منتهی به 1400/06"
$arr58[0,4] = "فصل سوم منتهی به 1400/09"
$arr58[0,5] = "فصل چهارم منتهی به 1400/12"
$arr58[0,6] = "فصل اول منتهی به 1401/03"
$arr58[0,7] = "فصل دوم منتهی به 1401/06"
$arr58[0,8] = "فصل سوم منتهی به 1401/09"
$arr58[0,9] = "فصل چهارم منتهی به 1401/12"
$ws.Range("E58:N58").Value = $arr58

$arr60 = New-Object 'object[,]' 1,10
$arr60[0,0] = 76419
$arr60[0,1] = 53923
$arr60[0,2] = 70108
$arr60[0,3] = 136585
$arr60[0,4] = 145086
$arr60[0,5] = "-"
$arr60[0,6] = 136456
$arr60[0,7] = 226121
$arr60[0,8] = 193485
$arr60[0,9] = 119473
$ws.Range("E60:N60").Value = $arr60

$arr61 = New-Object 'object[,]' 1,10
$arr61[0,0] = -6917
$arr61[0,1] = 15899
$arr61[0,2] = "-"
$arr61[0,3] = 579
$arr61[0,4] = 9112
$arr61[0,5] = "-"
$arr61[0,6] = 0
$arr61[0,7] = 16764
$arr61[0,8] = -16764
$arr61[0,9] = 27432
$ws.Range("E61:N61").Value = $arr61

$arr62 = New-Object 'object[,]' 1,10
$arr62[0,0] = 1575846
$arr62[0,1] = 1443556
$arr62[0,2] = 1363635
$arr62[0,3] = 1433088
$arr62[0,4] = 1791186
$arr62[0,5] = 1335131
$arr62[0,6] = "-"
$arr62[0,7] = "-"
$arr62[0,8] = "-"
$arr62[0,9] = "-"
$ws.Range("E62:N62").Value = $arr62

$arr63 = New-Object 'object[,]' 1,10
$arr63[0,0] = 32383
$arr63[0,1] = 6545
$arr63[0,2] = 23899
$arr63[0,3] = 112850
$arr63[0,4] = 3870
$arr63[0,5] = "-"
$arr63[0,6] = 21530
$arr63[0,7] = -118
$arr63[0,8] = 22090
$arr63[0,9] = -2181
$ws.Range("E63:N63").Value = $arr63

$arr64 = New-Object 'object[,]' 1,10
$arr64[0,0] = 0
$arr64[0,1] = 0
$arr64[0,2] = "-"
$arr64[0,3] = "-"
$arr64[0,4] = "-"
$arr64[0,5] = "-"
$arr64[0,6] = "-"
$arr64[0,7] = "-"
$arr64[0,8] = "-"
$arr64[0,9] = "-"
$ws.Range("E64:N64").Value = $arr64

$arr66 = New-Object 'object[,]' 1,10
$arr66[0,0] = 0
$arr66[0,1] = 0
$arr66[0,2] = "-"
$arr66[0,3] = 0
$arr66[0,4] = "-"
$arr66[0,5] = "-"
$arr66[0,6] = 0
$arr66[0,7] = 0
$arr66[0,8] = 0
$arr66[0,9] = 0
$ws.Range("E66:N66").Value = $arr66

$arr69 = New-Object 'object[,]' 1,10
$arr69[0,0] = 110589
$arr69[0,1] = 128636
$arr69[0,2] = 185478
$arr69[0,3] = 121522
$arr69[0,4] = 131791
$arr69[0,5] = "-"
$arr69[0,6] = 126791
$arr69[0,7] = 213892
$arr69[0,8] = 287378
$arr69[0,9] = 243946
$ws.Range("E69:N69").Value = $arr69

$arr70 = New-Object 'object[,]' 1,10
$arr70[0,0] = "-"
$arr70[0,1] = "-"
$arr70[0,2] = "-"
$arr70[0,3] = "-"
$arr70[0,4] = "-"
$arr70[0,5] = "-"
$arr70[0,6] = 1586731
$arr70[0,7] = 2587964
$arr70[0,8] = 3300472
$arr70[0,9] = 3263985
$ws.Range("E70:N70").Value = $arr70

$arr71 = New-Object 'object[,]' 1,10
$arr71[0,0] = 1788320
$arr71[0,1] = 1648559
$arr71[0,2] = 1643120
$arr71[0,3] = 1804624
$arr71[0,4] = 2081045
$arr71[0,5] = 1335131
$arr71[0,6] = 1871508
$arr71[0,7] = 3044623
$arr71[0,8] = 3786661
$arr71[0,9] = 3652655
$ws.Range("E71:N71").Value = $arr71

$arr75 = New-Object 'object[,]' 1,10
$arr75[0,0] = "فصل سوم منتهی به 1399/09"
$arr75[0,1] = "فصل چهارم منتهی به 1399/12"
$arr75[0,2] = "فصل اول منتهی به 1400/03"
$arr75[0,3] = "فصل دوم منتهی به 1400/06"
$arr75[0,4] = "فصل سوم منتهی به 1400/09"
$arr75[0,5] = "فصل چهارم منتهی به 1400/12"
$arr75[0,6] = "فصل اول منتهی به 1401/03"
$arr75[0,7] = "فصل دوم منتهی به 1401/06"
$arr75[0,8] = "فصل سوم منتهی به 1401/09"
$arr75[0,9] = "فصل چهارم منتهی به 1401/12"
$ws.Range("E75:N75").Value = $arr75

$arr77 = New-Object 'object[,]' 1,10
$arr77[0,0] = 3178034
$arr77[0,1] = 2145591
$arr77[0,2] = 2937568
$arr77[0,3] = 3691985
$arr77[0,4] = 3784293
$arr77[0,5] = 3846022
$arr77[0,6] = 3908458
$arr77[0,7] = 6062225
$arr77[0,8] = 6727105
$arr77[0,9] = 8091088
$ws.Range("E77:N77").Value = $arr77

$arr79 = New-Object 'object[,]' 1,10
$arr79[0,0] = 1798293
$arr79[0,1] = 1554128
$arr79[0,2] = 1595186
$arr79[0,3] = 1709631
$arr79[0,4] = 2030287
$arr79[0,5] = 1494421
$arr79[0,6] = "-"
$arr79[0,7] = "-"
$arr79[0,8] = "-"
$arr79[0,9] = "-"
$ws.Range("E79:N79").Value = $arr79

$arr80 = New-Object 'object[,]' 1,10
$arr80[0,0] = 161915000
$arr80[0,1] = 36564246
$arr80[0,2] = 165965278
$arr80[0,3] = 311740331
$arr80[0,4] = 10025907
$arr80[0,5] = -47905063
$arr80[0,6] = 60647887
$arr80[0,7] = -2000000
$arr80[0,8] = 380862069
$arr80[0,9] = 14938356
$ws.Range("E80:N80").Value = $arr80

$arr85 = New-Object 'object[,]' 1,10
$arr85[0,0] = 44288747
$arr85[0,1] = 60907197
$arr85[0,2] = 70955624
$arr85[0,3] = 56234151
$arr85[0,4] = 54890046
$arr85[0,5] = 73190679
$arr85[0,6] = 59359082
$arr85[0,7] = 136410714
$arr85[0,8] = 83273834
$arr85[0,9] = 112159080
$ws.Range("E85:N85").Value = $arr85

$arr86 = New-Object 'object[,]' 1,10
$arr86[0,0] = "-"
$arr86[0,1] = "-"
$arr86[0,2] = "-"
$arr86[0,3] = "-"
$arr86[0,4] = "-"
$arr86[0,5] = "-"
$arr86[0,6] = 1836790
$arr86[0,7] = 3411136
$arr86[0,8] = 2802439
$arr86[0,9] = 3982874
$ws.Range("E86:N86").Value = $arr86

$arr90 = New-Object 'object[,]' 1,10
$arr90[0,0] = "فصل سوم منتهی به 1399/09"
$arr90[0,1] = "فصل چهارم منتهی به 1399/12"
$arr90[0,2] = "فصل اول منتهی به 1400/03"
$arr90[0,3] = "فصل دوم منتهی به 1400/06"
$arr90[0,4] = "فصل سوم منتهی به 1400/09"
$arr90[0,5] = "فصل چهارم منتهی به 1400/12"
$arr90[0,6] = "فصل اول منتهی به 1401/03"
$arr90[0,7] = "فصل دوم منتهی به 1401/06"
$arr90[0,8] = "فصل سوم منتهی به 1401/09"
$arr90[0,9] = "فصل چهارم منتهی به 1401/12"
$ws.Range("E90:N90").Value = $arr90

$arr92 = New-Object 'object[,]' 1,10
$arr92[0,0] = -32075
$arr92[0,1] = -41159
$arr92[0,2] = -42324
$arr92[0,3] = -78740
$arr92[0,4] = -79667
$arr92[0,5] = -82251
$arr92[0,6] = -87557
$arr92[0,7] = -70358
$arr92[0,8] = -85053
$arr92[0,9] = -75311
$ws.Range("E92:N92").Value = $arr92

$arr93 = New-Object 'object[,]' 1,10
$arr93[0,0] = 10299
$arr93[0,1] = -10236
$arr93[0,2] = "-"
$arr93[0,3] = -1357
$arr93[0,4] = -7296
$arr93[0,5] = -5131
$arr93[0,6] = 0
$arr93[0,7] = -17540
$arr93[0,8] = 17540
$arr93[0,9] = -33116
$ws.Range("E93:N93").Value = $arr93

$arr94 = New-Object 'object[,]' 1,10
$arr94[0,0] = -580621
$arr94[0,1] = -648999
$arr94[0,2] = -637544
$arr94[0,3] = -914870
$arr94[0,4] = -894921
$arr94[0,5] = -954105
$arr94[0,6] = "-"
$arr94[0,7] = "-"
$arr94[0,8] = "-"
$arr94[0,9] = "-"
$ws.Range("E94:N94").Value = $arr94

$arr95 = New-Object 'object[,]' 1,10
$arr95[0,0] = -12274
$arr95[0,1] = -18168
$arr95[0,2] = -12435
$arr95[0,3] = -31122
$arr95[0,4] = -11468
$arr95[0,5] = 4656
$arr95[0,6] = -15284
$arr95[0,7] = 4172
$arr95[0,8] = -31377
$arr95[0,9] = 8676
$ws.Range("E95:N95").Value = $arr95

$arr96 = New-Object 'object[,]' 1,10
$arr96[0,0] = 0
$arr96[0,1] = 0
$arr96[0,2] = "-"
$arr96[0,3] = "-"
$arr96[0,4] = "-"
$arr96[0,5] = "-"
$arr96[0,6] = "-"
$arr96[0,7] = "-"
$arr96[0,8] = "-"
$arr96[0,9] = "-"
$ws.Range("E96:N96").Value = $arr96

$arr98 = New-Object 'object[,]' 1,10
$arr98[0,0] = 0
$arr98[0,1] = 0
$arr98[0,2] = "-"
$arr98[0,3] = 0
$arr98[0,4] = "-"
$arr98[0,5] = 0
$arr98[0,6] = 0
$arr98[0,7] = 0
$arr98[0,8] = 0
$arr98[0,9] = 0
$ws.Range("E98:N98").Value = $arr98

$arr100 = New-Object 'object[,]' 1,10
$arr100[0,0] = -60208
$arr100[0,1] = -73445
$arr100[0,2] = -85793
$arr100[0,3] = -94015
$arr100[0,4] = -104646
$arr100[0,5] = -151097
$arr100[0,6] = -106329
$arr100[0,7] = -108956
$arr100[0,8] = -193289
$arr100[0,9] = -200134
$ws.Range("E100:N100").Value = $arr100

$arr101 = New-Object 'object[,]' 1,10
$arr101[0,0] = "-"
$arr101[0,1] = "-"
$arr101[0,2] = "-"
$arr101[0,3] = "-"
$arr101[0,4] = "-"
$arr101[0,5] = "-"
$arr101[0,6] = -912500
$arr101[0,7] = -1032115
$arr101[0,8] = -1626126
$arr101[0,9] = -1721425
$ws.Range("E101:N101").Value = $arr101

$arr102 = New-Object 'object[,]' 1,10
$arr102[0,0] = -674879
$arr102[0,1] = -792007
$arr102[0,2] = -778096
$arr102[0,3] = -1120104
$arr102[0,4] = -1097998
$arr102[0,5] = -1187928
$arr102[0,6] = -1121670
$arr102[0,7] = -1224797
$arr102[0,8] = -1918305
$arr102[0,9] = -2021310
$ws.Range("E102:N102").Value = $arr102

$arr106 = New-Object 'object[,]' 1,10
$arr106[0,0] = "فصل سوم منتهی به 1399/09"
$arr106[0,1] = "فصل چهارم منتهی به 1399/12"
$arr106[0,2] = "فصل اول منتهی به 1400/03"
$arr106[0,3] = "فصل دوم منتهی به 1400/06"
$arr106[0,4] = "فصل سوم منتهی به 1400/09"
$arr106[0,5] = "فصل چهارم منتهی به 1400/12"
$arr106[0,6] = "فصل اول منتهی به 1401/03"
$arr106[0,7] = "فصل دوم منتهی به 1401/06"
$arr106[0,8] = "فصل سوم منتهی به 1401/09"
$arr106[0,9] = "فصل چهارم منتهی به 1401/12"
$ws.Range("E106:N106").Value = $arr106

$arr108 = New-Object 'object[,]' 1,10
$arr108[0,0] = 44344
$arr108[0,1] = 12764
$arr108[0,2] = 27784
$arr108[0,3] = 57845
$arr108[0,4] = 65419
$arr108[0,5] = 76882
$arr108[0,6] = 48899
$arr108[0,7] = 155763
$arr108[0,8] = 108432
$arr108[0,9] = 44162
$ws.Range("E108:N108").Value = $arr108

$arr109 = New-Object 'object[,]' 1,10
$arr109[0,0] = 3382
$arr109[0,1] = 5663
$arr109[0,2] = "-"
$arr109[0,3] = 349
$arr109[0,4] = 1816
$arr109[0,5] = 4596
$arr109[0,6] = 0
$arr109[0,7] = -776
$arr109[0,8] = 776
$arr109[0,9] = -5684
$ws.Range("E109:N109").Value = $arr109

$arr110 = New-Object 'object[,]' 1,10
$arr110[0,0] = 995225
$arr110[0,1] = 794557
$arr110[0,2] = 726091
$arr110[0,3] = 518218
$arr110[0,4] = 896265
$arr110[0,5] = 381026
$arr110[0,6] = "-"
$arr110[0,7] = "-"
$arr110[0,8] = "-"
$arr110[0,9] = "-"
$ws.Range("E110:N110").Value = $arr110

$arr111 = New-Object 'object[,]' 1,10
$arr111[0,0] = 20109
$arr111[0,1] = -11623
$arr111[0,2] = 11464
$arr111[0,3] = 81728
$arr111[0,4] = -7598
$arr111[0,5] = -2913
$arr111[0,6] = 6246
$arr111[0,7] = 4054
$arr111[0,8] = -9287
$arr111[0,9] = 6495
$ws.Range("E111:N111").Value = $arr111

$arr112 = New-Object 'object[,]' 1,10
$arr112[0,0] = 0
$arr112[0,1] = 0
$arr112[0,2] = "-"
$arr112[0,3] = "-"
$arr112[0,4] = "-"
$arr112[0,5] = "-"
$arr112[0,6] = "-"
$arr112[0,7] = "-"
$arr112[0,8] = "-"
$arr112[0,9] = "-"
$ws.Range("E112:N112").Value = $arr112

$arr114 = New-Object 'object[,]' 1,10
$arr114[0,0] = 0
$arr114[0,1] = 0
$arr114[0,2] = "-"
$arr114[0,3] = 0
$arr114[0,4] = "-"
$arr114[0,5] = 0
$arr114[0,6] = 0
$arr114[0,7] = 0
$arr114[0,8] = 0
$arr114[0,9] = 0
$ws.Range("E114:N114").Value = $arr114

$arr116 = New-Object 'object[,]' 1,10
$arr116[0,0] = 50381
$arr116[0,1] = 55191
$arr116[0,2] = 99685
$arr116[0,3] = 27507
$arr116[0,4] = 27145
$arr116[0,5] = 54642
$arr116[0,6] = 20462
$arr116[0,7] = 104936
$arr116[0,8] = 94089
$arr116[0,9] = 43812
$ws.Range("E116:N116").Value = $arr116

$arr117 = New-Object 'object[,]' 1,10
$arr117[0,0] = "-"
$arr117[0,1] = "-"
$arr117[0,2] = "-"
$arr117[0,3] = "-"
$arr117[0,4] = "-"
$arr117[0,5] = "-"
$arr117[0,6] = 674231
$arr117[0,7] = 1555849
$arr117[0,8] = 1674346
$arr117[0,9] = 1542560
$ws.Range("E117:N117").Value = $arr117

$arr118 = New-Object 'object[,]' 1,10
$arr118[0,0] = 1113441
$arr118[0,1] = 856552
$arr118[0,2] = 865024
$arr118[0,3] = 685647
$arr118[0,4] = 983047
$arr118[0,5] = 514233
$arr118[0,6] = 749838
$arr118[0,7] = 1819826
$arr118[0,8] = 1868356
$arr118[0,9] = 1631345
$ws.Range("E118:N118").Value = $arr118

